$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.514.56"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.575.96"
$ws.Range("E3").Value = "  -2.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.36"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.00"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.574.78"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.043.80"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.401.22"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.568.40"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.41"
$ws.Range("E19").Value = "  -6.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.75"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.04"
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  -9.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.707.36"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0989"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "529.15"
$ws.Range("E31").Value = "  -3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.12"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.63"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.12"
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("E46").Value = "  -4.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.35"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0763"
$ws.Range("E51").Value = "  -1.04%  "
